# Weekly refresh: prepend this week's two new price observations
# (Plátano, Feria Lagunitas de Puerto Montt) to the top of the data
# table, pushing the existing history down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh blank rows right above the current first data row
# (row 711); inserting twice at the same index pushes the previous
# insert down, opening up rows 711 and 712 for the new observations
# while the rest of the table (old 711..819) shifts down to 713..821.
$ws.Rows.Item(711).Insert()
$ws.Rows.Item(711).Insert()

# New row 711: "Sin especificar" / "Pintón" observation for 2023-04-18
$ws.Range("A711").Value = 4
$ws.Range("B711").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C711").Value = "Los Lagos"
$ws.Range("D711").Value = 45034
$ws.Range("E711").Value = 10
$ws.Range("F711").Value = "Fruta"
$ws.Range("G711").Value = 100108
$ws.Range("H711").Value = "Tropicales y subtropicales"
$ws.Range("I711").Value = 100108006
$ws.Range("J711").Value = "Plátano"
$ws.Range("K711").Value = "Sin especificar"
$ws.Range("L711").Value = "Pintón"
$ws.Range("M711").Value = 600
$ws.Range("N711").Value = 23000
$ws.Range("O711").Value = 23000
$ws.Range("P711").Value = 23000
$ws.Range("Q711").Value = "$/caja 20 kilos"
$ws.Range("R711").Value = "Ecuador"
$ws.Range("S711").Value = 1150
$ws.Range("T711").Value = 20

# New row 712: "Sin especificar" / "Primera Pintón" observation for 2023-04-18
$ws.Range("A712").Value = 4
$ws.Range("B712").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C712").Value = "Los Lagos"
$ws.Range("D712").Value = 45034
$ws.Range("E712").Value = 10
$ws.Range("F712").Value = "Fruta"
$ws.Range("G712").Value = 100108
$ws.Range("H712").Value = "Tropicales y subtropicales"
$ws.Range("I712").Value = 100108006
$ws.Range("J712").Value = "Plátano"
$ws.Range("K712").Value = "Sin especificar"
$ws.Range("L712").Value = "Primera Pintón"
$ws.Range("M712").Value = 1200
$ws.Range("N712").Value = 24000
$ws.Range("O712").Value = 25000
$ws.Range("P712").Value = 24500
$ws.Range("Q712").Value = "$/caja 20 kilos"
$ws.Range("R712").Value = "Ecuador"
$ws.Range("S712").Value = 1225
$ws.Range("T712").Value = 20
